$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name
$ws.Name = "Through 2022-09-03"

# Update header label for September row
$ws.Range("A10").Value = "September (through 09-03)"

# Update September row (row 10) values
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 8
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = 19

# Update Total row (row 11) values
$ws.Range("B11").Value = 195
$ws.Range("C11").Value = 385
$ws.Range("D11").Value = 561
$ws.Range("E11").Value = 495
$ws.Range("F11").Value = 360
$ws.Range("G11").Value = 792
$ws.Range("H11").Value = 1080
$ws.Range("I11").Value = 1158
